# correct inputs, add AR28
# Kevin Cahill's personnel role is corrected from "technician" to "creator".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Personnel")
$ws.Activate()

# Reflect the cell the edit was made on (matches the saved selection state).
[void]$ws.Range("G7").Select()

# The actual data correction: row 6 (Kevin Cahill) role -> creator
$ws.Range("G6").Value = "creator"
